$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.287853240966797
$ws.Range("B1").Value = 2.995970964431763
$ws.Range("C1").Value = 1.461484789848328
$ws.Range("D1").Value = 1.197767972946167
$ws.Range("E1").Value = 1.260863423347473
